# Update the "想去人数" (want-to-go count) column F on each sheet to the
# refreshed values captured at commit 456a3b4 (gh-pages data regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1787
$ws1.Range("F8").Value  = 1753
$ws1.Range("F18").Value = 129
$ws1.Range("F19").Value = 448
$ws1.Range("F20").Value = 448
$ws1.Range("F22").Value = 375
$ws1.Range("F24").Value = 1140
$ws1.Range("F26").Value = 1252
$ws1.Range("F28").Value = 1443
$ws1.Range("F35").Value = 750
$ws1.Range("F38").Value = 0
$ws1.Range("F41").Value = 1283
$ws1.Range("F45").Value = 58

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 74
$ws2.Range("F20").Value = 2326
$ws2.Range("F24").Value = 111

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value  = 2592
$ws3.Range("F7").Value  = 4290
$ws3.Range("F11").Value = 416
$ws3.Range("F12").Value = 307
$ws3.Range("F13").Value = 256
$ws3.Range("F14").Value = 97

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1787
$ws4.Range("F6").Value  = 4290
$ws4.Range("F8").Value  = 416
$ws4.Range("F9").Value  = 1753
$ws4.Range("F12").Value = 256
$ws4.Range("F13").Value = 97
$ws4.Range("F18").Value = 448
$ws4.Range("F19").Value = 448
$ws4.Range("F21").Value = 375
$ws4.Range("F22").Value = 2326
$ws4.Range("F24").Value = 1140
$ws4.Range("F26").Value = 1252
$ws4.Range("F27").Value = 111
$ws4.Range("F29").Value = 1443
$ws4.Range("F34").Value = 750
$ws4.Range("F44").Value = 58
